# ----------------------------------------------------------------------------
# Edit script: initial cleaning and exploration
#
# - Corrects ~50 mis-keyed values in column G ("WORRIED?") for rows 159-218
# - Changes the font of A158 ("1703SM") and G158 to Arial
#   (this causes Excel to register two new font/style entries and bumps the
#   row height of row 158 from the default 14.25 to 18.0, matching a
#   user font change made while cleaning the data)
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply font change to A158 and G158 (Arial), which Excel reflects as new
# style/font entries, and grow row 158 to fit the new font (18pt).
$ws.Range("A158").Font.Name = "Arial"
$ws.Range("G158").Font.Name = "Arial"
$ws.Rows.Item(158).RowHeight = 18.0

# Correct the mis-entered "WORRIED?" (column G) values for rows 159-218.
$ws.Range("G159").Value = 4.0
$ws.Range("G160").Value = 5.0
$ws.Range("G161").Value = 4.0
$ws.Range("G162").Value = 1.0
$ws.Range("G163").Value = 2.0
$ws.Range("G164").Value = 4.0
$ws.Range("G166").Value = 2.0
$ws.Range("G167").Value = 5.0
$ws.Range("G169").Value = 4.0
$ws.Range("G170").Value = 1.0
$ws.Range("G171").Value = 5.0
$ws.Range("G173").Value = 4.0
$ws.Range("G174").Value = 5.0
$ws.Range("G175").Value = 4.0
$ws.Range("G176").Value = 1.0
$ws.Range("G177").Value = 4.0
$ws.Range("G178").Value = 2.0
$ws.Range("G179").Value = 4.0
$ws.Range("G180").Value = 2.0
$ws.Range("G181").Value = 4.0
$ws.Range("G182").Value = 5.0
$ws.Range("G185").Value = 4.0
$ws.Range("G186").Value = 5.0
$ws.Range("G187").Value = 4.0
$ws.Range("G188").Value = 2.0
$ws.Range("G189").Value = 4.0
$ws.Range("G191").Value = 5.0
$ws.Range("G192").Value = 2.0
$ws.Range("G194").Value = 5.0
$ws.Range("G195").Value = 4.0
$ws.Range("G196").Value = 2.0
$ws.Range("G197").Value = 4.0
$ws.Range("G198").Value = 1.0
$ws.Range("G199").Value = 4.0
$ws.Range("G200").Value = 3.0
$ws.Range("G201").Value = 4.0
$ws.Range("G202").Value = 3.0
$ws.Range("G204").Value = 4.0
$ws.Range("G205").Value = 3.0
$ws.Range("G206").Value = 4.0
$ws.Range("G207").Value = 1.0
$ws.Range("G209").Value = 4.0
$ws.Range("G210").Value = 3.0
$ws.Range("G212").Value = 2.0
$ws.Range("G213").Value = 4.0
$ws.Range("G214").Value = 5.0
$ws.Range("G215").Value = 1.0
$ws.Range("G216").Value = 4.0
$ws.Range("G217").Value = 2.0
$ws.Range("G218").Value = 1.0
